$wb = $excel.ActiveWorkbook

# Create Sheet3 by copying Sheet2 (inherits sheetPr/outline, column widths,
# page margins and the 4 header merges) and placing it right after Sheet2.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Copy($null, $ws2)
$sheet3 = $wb.Worksheets.Item(3)
$sheet3.Name = "Sheet3"

# Create Sheet4 the same way, by copying the freshly made Sheet3, placed after it.
$sheet3.Copy($null, $sheet3)
$sheet4 = $wb.Worksheets.Item(4)
$sheet4.Name = "Sheet4"

# Overwrite Sheet3 data (labels identical to Sheet2's layout; values are new).
$sheet3.Cells.Item(1, 1).Value = "Кількість елементів"
$sheet3.Cells.Item(1, 2).Value = "Середній час"
$sheet3.Cells.Item(1, 6).Value = "Кількість елементів"
$sheet3.Cells.Item(1, 7).Value = "Точність"
$sheet3.Cells.Item(2, 2).Value = "B&B"
$sheet3.Cells.Item(2, 3).Value = "Greed"
$sheet3.Cells.Item(2, 4).Value = "ACO"
$sheet3.Cells.Item(2, 7).Value = "B&B"
$sheet3.Cells.Item(2, 8).Value = "Greed"
$sheet3.Cells.Item(2, 9).Value = "ACO"
$sheet3.Cells.Item(3, 1).Value = 5
$sheet3.Cells.Item(3, 2).Value = [double]"3.100000321865082e-05"
$sheet3.Cells.Item(3, 3).Value = 0.002045609999913722
$sheet3.Cells.Item(3, 4).Value = 0.08273572000180138
$sheet3.Cells.Item(3, 6).Value = 5
$sheet3.Cells.Item(3, 7).Value = 1.2
$sheet3.Cells.Item(3, 8).Value = 0
$sheet3.Cells.Item(3, 9).Value = 0
$sheet3.Cells.Item(4, 1).Value = 6
$sheet3.Cells.Item(4, 2).Value = [double]"3.619999915827066e-05"
$sheet3.Cells.Item(4, 3).Value = 0.003190910002740566
$sheet3.Cells.Item(4, 4).Value = 0.1229132299980847
$sheet3.Cells.Item(4, 6).Value = 6
$sheet3.Cells.Item(4, 7).Value = 1.4
$sheet3.Cells.Item(4, 8).Value = 0
$sheet3.Cells.Item(4, 9).Value = 0
$sheet3.Cells.Item(5, 1).Value = 7
$sheet3.Cells.Item(5, 2).Value = [double]"4.355000128271058e-05"
$sheet3.Cells.Item(5, 3).Value = 0.01396054999931948
$sheet3.Cells.Item(5, 4).Value = 0.1651210500000161
$sheet3.Cells.Item(5, 6).Value = 7
$sheet3.Cells.Item(5, 7).Value = 1.5
$sheet3.Cells.Item(5, 8).Value = 0
$sheet3.Cells.Item(5, 9).Value = 0
$sheet3.Cells.Item(6, 1).Value = 8
$sheet3.Cells.Item(6, 2).Value = [double]"5.860000237589702e-05"
$sheet3.Cells.Item(6, 3).Value = 0.03141236000083154
$sheet3.Cells.Item(6, 4).Value = 0.1737159000040265
$sheet3.Cells.Item(6, 6).Value = 8
$sheet3.Cells.Item(6, 7).Value = 2
$sheet3.Cells.Item(6, 8).Value = 0
$sheet3.Cells.Item(6, 9).Value = 0.2
$sheet3.Cells.Item(7, 1).Value = 9
$sheet3.Cells.Item(7, 2).Value = [double]"6.352000054903329e-05"
$sheet3.Cells.Item(7, 3).Value = 0.05206127000274137
$sheet3.Cells.Item(7, 4).Value = 0.2363194600024144
$sheet3.Cells.Item(7, 6).Value = 9
$sheet3.Cells.Item(7, 7).Value = 2.7
$sheet3.Cells.Item(7, 8).Value = 0
$sheet3.Cells.Item(7, 9).Value = 0.3
$sheet3.Cells.Item(8, 1).Value = 10
$sheet3.Cells.Item(8, 2).Value = [double]"8.221999887609854e-05"
$sheet3.Cells.Item(8, 3).Value = 0.6565208600033656
$sheet3.Cells.Item(8, 4).Value = 0.2633120999991661
$sheet3.Cells.Item(8, 6).Value = 10
$sheet3.Cells.Item(8, 7).Value = 2
$sheet3.Cells.Item(8, 8).Value = 0
$sheet3.Cells.Item(8, 9).Value = 0.4
$sheet3.Cells.Item(9, 1).Value = 11
$sheet3.Cells.Item(9, 2).Value = 0.0001598099988768809
$sheet3.Cells.Item(9, 3).Value = 7.037033379996137
$sheet3.Cells.Item(9, 4).Value = 0.8054246999978204
$sheet3.Cells.Item(9, 6).Value = 11
$sheet3.Cells.Item(9, 7).Value = 1.7
$sheet3.Cells.Item(9, 8).Value = 0
$sheet3.Cells.Item(9, 9).Value = 0.4
$sheet3.Cells.Item(10, 1).Value = 12
$sheet3.Cells.Item(10, 2).Value = 0.0001328599959379062
$sheet3.Cells.Item(10, 3).Value = 7.806809709999652
$sheet3.Cells.Item(10, 4).Value = 0.4187305200000992
$sheet3.Cells.Item(10, 6).Value = 12
$sheet3.Cells.Item(10, 7).Value = 2.2
$sheet3.Cells.Item(10, 8).Value = 0
$sheet3.Cells.Item(10, 9).Value = 0.5
$sheet3.Cells.Item(11, 1).Value = 13
$sheet3.Cells.Item(11, 2).Value = 0.0001266500024939887
$sheet3.Cells.Item(11, 3).Value = 116.5706978299975
$sheet3.Cells.Item(11, 4).Value = 0.4791568699976779
$sheet3.Cells.Item(11, 6).Value = 13
$sheet3.Cells.Item(11, 7).Value = 3.1
$sheet3.Cells.Item(11, 8).Value = 0
$sheet3.Cells.Item(11, 9).Value = 0.8

# Overwrite Sheet4 data.
$sheet4.Cells.Item(1, 1).Value = "Кількість елементів"
$sheet4.Cells.Item(1, 2).Value = "Середній час"
$sheet4.Cells.Item(1, 6).Value = "Кількість елементів"
$sheet4.Cells.Item(1, 7).Value = "Точність"
$sheet4.Cells.Item(2, 2).Value = "B&B"
$sheet4.Cells.Item(2, 3).Value = "Greed"
$sheet4.Cells.Item(2, 4).Value = "ACO"
$sheet4.Cells.Item(2, 7).Value = "B&B"
$sheet4.Cells.Item(2, 8).Value = "Greed"
$sheet4.Cells.Item(2, 9).Value = "ACO"
$sheet4.Cells.Item(3, 1).Value = 5
$sheet4.Cells.Item(3, 2).Value = [double]"3.027000639121979e-05"
$sheet4.Cells.Item(3, 3).Value = 0.001713790002395399
$sheet4.Cells.Item(3, 4).Value = 0.08025141999823973
$sheet4.Cells.Item(3, 6).Value = 5
$sheet4.Cells.Item(3, 7).Value = 0.7
$sheet4.Cells.Item(3, 8).Value = 0
$sheet4.Cells.Item(3, 9).Value = 0
$sheet4.Cells.Item(4, 1).Value = 7
$sheet4.Cells.Item(4, 2).Value = [double]"5.331000429578126e-05"
$sheet4.Cells.Item(4, 3).Value = 0.01992319999844767
$sheet4.Cells.Item(4, 4).Value = 0.1472878600063268
$sheet4.Cells.Item(4, 6).Value = 7
$sheet4.Cells.Item(4, 7).Value = 1.9
$sheet4.Cells.Item(4, 8).Value = 0
$sheet4.Cells.Item(4, 9).Value = 0.1
$sheet4.Cells.Item(5, 1).Value = 9
$sheet4.Cells.Item(5, 2).Value = [double]"7.309000066015869e-05"
$sheet4.Cells.Item(5, 3).Value = 0.1057229499972891
$sheet4.Cells.Item(5, 4).Value = 0.2366165999992518
$sheet4.Cells.Item(5, 6).Value = 9
$sheet4.Cells.Item(5, 7).Value = 1.6
$sheet4.Cells.Item(5, 8).Value = 0
$sheet4.Cells.Item(5, 9).Value = 0.3
$sheet4.Cells.Item(6, 1).Value = 11
$sheet4.Cells.Item(6, 2).Value = [double]"9.344999853055924e-05"
$sheet4.Cells.Item(6, 3).Value = 41.53612096000579
$sheet4.Cells.Item(6, 4).Value = 0.3472341199987568
$sheet4.Cells.Item(6, 6).Value = 11
$sheet4.Cells.Item(6, 7).Value = 1.7
$sheet4.Cells.Item(6, 8).Value = 0
$sheet4.Cells.Item(6, 9).Value = 0.4
$sheet4.Cells.Item(7, 1).Value = 13
$sheet4.Cells.Item(7, 2).Value = 0.0001149299961980432
$sheet4.Cells.Item(7, 3).Value = 108.860699790009
$sheet4.Cells.Item(7, 4).Value = 0.651249509997433
$sheet4.Cells.Item(7, 6).Value = 13
$sheet4.Cells.Item(7, 7).Value = 3.4
$sheet4.Cells.Item(7, 8).Value = 0
$sheet4.Cells.Item(7, 9).Value = 1.1

# Restore Sheet1 as the active/selected sheet (unchanged from the original file).
$wb.Worksheets.Item(1).Activate()
